# Update column C (end date) from 45174 (2023-09-05) to 45175 (2023-09-06)
# for all data rows on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 32) {
    $lastRow = 32
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
